# Working version with docstring and params for filenames
#
# Reproduces, via Excel COM-interop, the edits captured in the target diff:
#   1. Resize/reposition the workbook window (bookViews/workbookView).
#   2. Sheet1!L18's literal value changes from 21.543 to 41.543 -- every
#      dependent M:R formula on rows 18-27 recalculates automatically from
#      this single edit.
#   3. Sheet1!L19 picks up a new cell style: a thin border on its right edge.
#   4. The active sheet's view/selection moves from N14 (topLeftCell F1) to
#      L19 (topLeftCell I15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- 1. Resize/reposition the workbook's window -------------------------
$win = $excel.ActiveWindow
$win.Left = 15000
$win.Top = 460
$win.Width = 18420
$win.Height = 8860

# --- 2. Update the measured value in L18 ---------------------------------
# (J18-I18)/L18*100*10^3 and friends recalc automatically; this is the
# single edit that drives every M18:R27 number in the diff.
$ws.Range("L18").Value = 41.543

# --- 3. Give L19 a thin right border (new cell style) --------------------
$ws.Range("L19").Borders.Item(10).LineStyle = 1
$ws.Range("L19").Borders.Item(10).Weight = 2

# --- 4. Scroll the view and move the active selection to L19 -------------
$win.ScrollRow = 15
$win.ScrollColumn = 9
[void]$ws.Range("L19").Select()

Write-Host "edit complete"
